$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("data")

# --- 1. Refresh the per-row query timestamps on the "data" sheet (column F) ---
$times = @(
    "2021-10-05 14:34:15.009248",
    "2021-10-05 14:34:15.009258",
    "2021-10-05 14:34:15.009262",
    "2021-10-05 14:34:15.009265",
    "2021-10-05 14:34:15.009268",
    "2021-10-05 14:34:15.009271",
    "2021-10-05 14:34:15.009273",
    "2021-10-05 14:34:15.009276",
    "2021-10-05 14:34:15.009279",
    "2021-10-05 14:34:15.009282",
    "2021-10-05 14:34:15.009285",
    "2021-10-05 14:34:15.009288",
    "2021-10-05 14:34:15.009290",
    "2021-10-05 14:34:15.009293",
    "2021-10-05 14:34:15.009296",
    "2021-10-05 14:34:15.009298",
    "2021-10-05 14:34:15.009301",
    "2021-10-05 14:34:15.009304",
    "2021-10-05 14:34:15.009307",
    "2021-10-05 14:34:15.009310",
    "2021-10-05 14:34:15.009313",
    "2021-10-05 14:34:15.009316",
    "2021-10-05 14:34:15.009319",
    "2021-10-05 14:34:15.009322",
    "2021-10-05 14:34:15.009325",
    "2021-10-05 14:34:15.009328",
    "2021-10-05 14:34:15.009330",
    "2021-10-05 14:34:15.009333"
)

$row = 2
foreach ($t in $times) {
    $dataSheet.Cells.Item($row, 6).Value = $t
    $row = $row + 1
}

# --- 2. Add a new "metadata" tab (placed after "data") describing the query ---
$ws = $wb.Worksheets.Add($null, $dataSheet)
$ws.Name = "metadata"

# Copy the header look & feel (bold, bordered, centered) from the "data" sheet
# so the new tab's style reuses the existing style definition instead of
# fabricating a new one.
$dataSheet.Range("B1:F1").Copy()
$ws.Range("B1:F1").PasteSpecial(-4122)
$dataSheet.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)
$dataSheet.Range("A2").Copy()
$ws.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("B1").Value = "data_name"
$ws.Range("C1").Value = "data_id"
$ws.Range("D1").Value = "data_version"
$ws.Range("E1").Value = "data_version_created"
$ws.Range("F1").Value = "panel_query_time"
$ws.Range("G1").Value = "panel_get_request"

$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "Iron metabolism disorders"
$ws.Range("C2").Value = 3469
$ws.Range("D2").Value = "'0.27"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "2021-09-14T05:56:25.467619Z"
$ws.Range("F2").Value = "2021-10-05 14:34:15.005649"
$ws.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/3469/?format=json"

# Restore "data" as the active sheet (matches original selection state).
$dataSheet.Activate()
